$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New student record (row 4)
$ws.Range("A4").Value = 206068
$ws.Range("B4").Value = "Omoniyi"
$ws.Range("C4").Value = "omoniyi@gmail.com"
[void]$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:omoniyi@gmail.com")
$ws.Range("C4").Style = $ws.Range("C3").Style
$ws.Range("D4").Value = "adult edu"
$ws.Range("E4").Value = 300
$ws.Range("F4").Value = "eiueh485"

# Match the author's final selection state
[void]$ws.Range("A11").Select()
